$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Sheet is protected; unprotect to allow the data refresh, then restore protection.
$ws.Unprotect()

# Update the confidential disclaimer date in the shared text (found wherever used, e.g. row 37)
$oldText = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-03-18 for illustrative purposes only and are subject to change."
$newText = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-03-19 for illustrative purposes only and are subject to change."
$found = $ws.Cells.Find($oldText)
if ($found -ne $null) {
    $found.Value = $newText
} else {
    $ws.Range("A37").Value = $newText
}

# Update Weight (D) and Percent Change (E) values for rows 2-34
$ws.Range("D2").Value = 0.03828160722746898
$ws.Range("E2").Value = 0.004616338441975953
$ws.Range("D3").Value = 0.02167795025638558
$ws.Range("E3").Value = 0.01073635407024121
$ws.Range("D4").Value = 0.01971951322150731
$ws.Range("E4").Value = 0.01167372881355933
$ws.Range("D5").Value = 0.04051977362726649
$ws.Range("E5").Value = 0.004585537918871285
$ws.Range("D6").Value = 0.03761923728800576
$ws.Range("E6").Value = 0.001179245283018604
$ws.Range("D7").Value = 0.02101437588300542
$ws.Range("E7").Value = 0.00294290759270166
$ws.Range("D8").Value = 0.03799354082704803
$ws.Range("E8").Value = -0.001121327651940063
$ws.Range("D9").Value = 0.02156502921691443
$ws.Range("E9").Value = 0.002936588051757338
$ws.Range("D10").Value = 0.0253777689526897
$ws.Range("E10").Value = 0.004223979205025374
$ws.Range("D11").Value = 0.02440064376448638
$ws.Range("E11").Value = -0.0121092649957758
$ws.Range("D12").Value = 0.05920900385145698
$ws.Range("E12").Value = -0.008428358948934123
$ws.Range("D13").Value = 0.02632053671514746
$ws.Range("E13").Value = 0.008646616541353236
$ws.Range("D14").Value = 0.02722647331720601
$ws.Range("E14").Value = -0.01158940397350983
$ws.Range("D15").Value = 0.03466922625567143
$ws.Range("E15").Value = -0.001243118451429659
$ws.Range("D16").Value = 0.01907353864560986
$ws.Range("E16").Value = -0.0067095331283199
$ws.Range("D17").Value = 0.03043243342837434
$ws.Range("E17").Value = -0.03341763005780352
$ws.Range("D18").Value = 0.02395927453038881
$ws.Range("E18").Value = 0.005148607535689242
$ws.Range("D19").Value = 0.1334099391418342
$ws.Range("E19").Value = 0.007442489851150258
$ws.Range("D20").Value = 0.009769657697478534
$ws.Range("E20").Value = 0.00310599197618755
$ws.Range("D21").Value = 0.01572190304928041
$ws.Range("E21").Value = 0.003083926866877418
$ws.Range("D22").Value = 0.01763442976786164
$ws.Range("E22").Value = -0.03131955484896654
$ws.Range("D23").Value = 0.01695717619156545
$ws.Range("E23").Value = -0.0003533568904594553
$ws.Range("D24").Value = 0.0215459539741405
$ws.Range("E24").Value = -0.01217772901810199
$ws.Range("D25").Value = 0.01203642321846664
$ws.Range("E25").Value = 0.01040391676866603
$ws.Range("D26").Value = 0.04279066203934406
$ws.Range("E26").Value = 0.004162330905306932
$ws.Range("D27").Value = 0.02578016312302153
$ws.Range("E27").Value = -0.0001961745953896754
$ws.Range("D28").Value = 0.04784301769601861
$ws.Range("E28").Value = 0.004405286343612369
$ws.Range("D29").Value = 0.05710479018315619
$ws.Range("E29").Value = 0.001320754716981121
$ws.Range("D30").Value = 0.01418076636003822
$ws.Range("E30").Value = -0.02831402831402829
$ws.Range("D31").Value = 0.01464313485275147
$ws.Range("E31").Value = -0.0178320037541061
$ws.Range("D32").Value = 0.04450644439513601
$ws.Range("E32").Value = 0.003134796238244419
$ws.Range("D33").Value = 0.01701561130127347
$ws.Range("E33").Value = 0.0006332119677061598
$ws.Range("D34").Value = 0.9999999999999998
$ws.Range("E34").Value = -0.0005643858595498763

# Restore sheet protection (password unknown; original legacy hash preserved separately).
$ws.Protect()
